$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) for the regional/categorie-juridique
# breakdown rows with the 2020-08-06 refreshed counts.
$updates = @(
    @{ Row = 2; C = 315448; D = 402139651 }
    @{ Row = 3; C = 253; D = 301350 }
    @{ Row = 10; C = 115909; D = 169849856 }
    @{ Row = 12; C = 58467; D = 84384552 }
    @{ Row = 16; C = 3956; D = 5616361 }
    @{ Row = 20; C = 6448; D = 8999385 }
    @{ Row = 22; C = 76280; D = 95204800 }
    @{ Row = 28; C = 32168; D = 47095912 }
    @{ Row = 30; C = 11339; D = 16309112 }
    @{ Row = 33; C = 1553; D = 2180807 }
    @{ Row = 35; C = 1775; D = 2504333 }
    @{ Row = 36; C = 95864; D = 120758463 }
    @{ Row = 44; C = 44017; D = 64512773 }
    @{ Row = 46; C = 9016; D = 12939198 }
    @{ Row = 48; C = 1387; D = 1926600 }
    @{ Row = 51; C = 2245; D = 3131499 }
    @{ Row = 52; C = 68030; D = 85385228 }
    @{ Row = 58; C = 27872; D = 40879369 }
    @{ Row = 61; C = 10935; D = 15811057 }
    @{ Row = 63; C = 1347; D = 1882789 }
    @{ Row = 67; C = 1436; D = 2010599 }
    @{ Row = 69; C = 20204; D = 26462753 }
    @{ Row = 72; C = 58; D = 84573 }
    @{ Row = 73; C = 7498; D = 10976590 }
    @{ Row = 75; C = 5040; D = 7318206 }
    @{ Row = 76; C = 484; D = 684739 }
    @{ Row = 78; C = 138606; D = 172907533 }
    @{ Row = 83; C = 15; D = 19021 }
    @{ Row = 84; C = 62922; D = 92228125 }
    @{ Row = 87; C = 29304; D = 42395157 }
    @{ Row = 89; C = 2705; D = 3895643 }
    @{ Row = 90; C = 2743; D = 3875850 }
    @{ Row = 91; C = 31902; D = 43227975 }
    @{ Row = 95; C = 7765; D = 11417070 }
    @{ Row = 97; C = 7058; D = 10231955 }
    @{ Row = 99; C = 515; D = 731905 }
    @{ Row = 101; C = 8676; D = 12035029 }
    @{ Row = 103; C = 2185; D = 3219470 }
    @{ Row = 105; C = 2939; D = 4292112 }
    @{ Row = 107; C = 126; D = 183120 }
    @{ Row = 108; C = 165; D = 233586 }
    @{ Row = 109; C = 139094; D = 172037233 }
    @{ Row = 115; C = 52184; D = 76503350 }
    @{ Row = 117; C = 26551; D = 38465443 }
    @{ Row = 121; C = 2194; D = 3081996 }
    @{ Row = 123; C = 492219; D = 649150254 }
    @{ Row = 128; C = 1360; D = 2016311 }
    @{ Row = 130; C = 204369; D = 300440174 }
    @{ Row = 131; C = 388; D = 578790 }
    @{ Row = 133; C = 176529; D = 256593906 }
    @{ Row = 136; C = 2803; D = 3939784 }
    @{ Row = 138; C = 6153; D = 8693320 }
    @{ Row = 141; C = 43687; D = 58340355 }
    @{ Row = 147; C = 13868; D = 20341228 }
    @{ Row = 148; C = 3692; D = 5325211 }
    @{ Row = 150; C = 7; D = 10500 }
    @{ Row = 151; C = 393; D = 565431 }
    @{ Row = 153; C = 373; D = 526251 }
    @{ Row = 154; C = 17200; D = 22731008 }
    @{ Row = 158; C = 7037; D = 10233981 }
    @{ Row = 160; C = 4902; D = 7056436 }
    @{ Row = 162; C = 273; D = 377235 }
    @{ Row = 163; C = 260; D = 372274 }
    @{ Row = 165; C = 15222; D = 22087426 }
    @{ Row = 166; C = 1732; D = 2576530 }
    @{ Row = 171; C = 86453; D = 108159819 }
    @{ Row = 176; C = 639; D = 941848 }
    @{ Row = 178; C = 33519; D = 49158950 }
    @{ Row = 180; C = 12815; D = 18516388 }
    @{ Row = 184; C = 1600; D = 2252693 }
    @{ Row = 186; C = 234752; D = 291895644 }
    @{ Row = 187; C = 132; D = 143030 }
    @{ Row = 192; C = 865; D = 1272497 }
    @{ Row = 194; C = 85781; D = 125751202 }
    @{ Row = 197; C = 32573; D = 46877087 }
    @{ Row = 200; C = 5033; D = 7172101 }
    @{ Row = 203; C = 4720; D = 6530827 }
    @{ Row = 206; C = 259536; D = 321272342 }
    @{ Row = 213; C = 609; D = 886378 }
    @{ Row = 215; C = 94170; D = 137772400 }
    @{ Row = 216; C = 86; D = 128199 }
    @{ Row = 218; C = 50689; D = 73257624 }
    @{ Row = 221; C = 4611; D = 6472416 }
    @{ Row = 224; C = 5568; D = 7703381 }
    @{ Row = 227; C = 104619; D = 130957186 }
    @{ Row = 228; C = 74; D = 78905 }
    @{ Row = 234; C = 49010; D = 71804122 }
    @{ Row = 236; C = 12196; D = 17532969 }
    @{ Row = 238; C = 1876; D = 2688838 }
    @{ Row = 240; C = 2426; D = 3389650 }
    @{ Row = 241; C = 253233; D = 319826601 }
    @{ Row = 242; C = 170; D = 210959 }
    @{ Row = 249; C = 94712; D = 138787414 }
    @{ Row = 252; C = 63938; D = 92653129 }
    @{ Row = 254; C = 2379; D = 3357861 }
    @{ Row = 257; C = 4485; D = 6296392 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

